$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.227.51"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "2.620.59"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.90"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.68"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.598"
$ws.Range("E7").Value = "  +1.25%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +3.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.18"
$ws.Range("E10").Value = "  +2.86%  "
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.23"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").Value = "3.021.62"
$ws.Range("E15").Value = "  +2.16%  "
$ws.Range("D16").Value = "2.616.78"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.921"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "46.497.49"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  +3.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.84"
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "294.28"
$ws.Range("E23").Value = "  +16.84%  "
$ws.Range("E24").Value = "  +4.32%  "
$ws.Range("E25").Value = "  +3.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.89"
$ws.Range("E27").Value = "  +8.76%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.85"
$ws.Range("E30").Value = "  +5.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "39.32"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("E32").Value = "  -1.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.30"
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.67"
$ws.Range("E35").Value = "  +4.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.25"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +2.64%  "
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("E39").Value = "  +7.04%  "
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.82"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("E43").Value = "  +2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.67"
$ws.Range("E44").Value = "  +11.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.05"
$ws.Range("E45").Value = "  -1.90%  "
$ws.Range("D46").Value = "2.132.06"
$ws.Range("E46").Value = "  +4.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.48"
$ws.Range("E47").Value = "  +7.63%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.52"
$ws.Range("E49").Value = "  +3.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "110.16"
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.202"
$ws.Range("E51").Value = "  +2.33%  "
